# Update 16 10 2021
# Applies the "uang kas kecil - asrama" Oktober rekap update:
#  - correct two dates (12 Oct -> 13 Oct) for entries 15 & 16
#  - add six new kas-kecil transactions (entries 17-22, 13-16 Oct)
#  - clear out the now-unused trailing shared-formula rows (36-41)
#  - update the saved view (scroll position / zoom / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the two mis-dated rows (entries 15 & 16) ---------------------------
$ws.Range("B23").Value = 44482
$ws.Range("B24").Value = 44482

# --- new transaction rows (17-22) ------------------------------------------
# row 25 - entry 17
$ws.Range("A25").Value = 17
$ws.Range("B25").Value = 44483
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 84500
$ws.Range("F25").Value = "uang beli tespan, voltmeter sama lem G"
$ws.Range("G25").Value = "Saferius sama Hosea"

# row 26 - entry 18
$ws.Range("A26").Value = 18
$ws.Range("B26").Value = 44485
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 32000
$ws.Range("F26").Value = "uang beli tabung gas 3 kg (2 tabung)"
$ws.Range("G26").Value = "jhonan sama peter"

# row 27 - entry 19
$ws.Range("A27").Value = 19
$ws.Range("B27").Value = 44485
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 31700
$ws.Range("F27").Value = "uang beli minyak 2L"
$ws.Range("G27").Value = "jhonan sama peter"

# row 28 - entry 20
$ws.Range("A28").Value = 20
$ws.Range("B28").Value = 44485
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 55000
$ws.Range("F28").Value = "uang beli DHS machine head (puteran gitar)"
$ws.Range("G28").Value = "Saferius sama tondo"

# row 29 - entry 21
$ws.Range("A29").Value = 21
$ws.Range("B29").Value = 44485
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 10000
$ws.Range("F29").Value = "uang beli bensin pertalite"
$ws.Range("G29").Value = "Saferius sama tondo"

# row 30 - entry 22
$ws.Range("A30").Value = 22
$ws.Range("B30").Value = 44485
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 2000
$ws.Range("F30").Value = "uang parkir"
$ws.Range("G30").Value = "Saferius sama tondo"

# --- rows 36-41 no longer carry the running-total formula -------------------
$ws.Range("E36:E41").ClearContents()

# --- saved view: scroll/zoom/selection update -------------------------------
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D31").Select()
$excel.ActiveWindow.Zoom = 90
